$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AutoFilter over the original table range (applied before the new
#     rows are appended so the filter range doesn't auto-expand) ---
$ws.Range("A1:H29").AutoFilter()

# --- New rows 30-34: "CLDRJar" interactive scenario test cases ---
# Values are written in the same order the shared-string table records
# them so the rebuilt workbook's sst ordering matches the authored file.

$ws.Range("H30").Value = "CLDRJar"
$ws.Range("B30").Value = "Check the number data can be fetched"
$ws.Range("C30").Value = "number"
$ws.Range("G32").Value = 'Il y a 100 000 fichiers sur "MyDisk".'
$ws.Range("B31").Value = "Check the percent data can be fetched"
$ws.Range("B32").Value = "Check the plurals data can be fetched"
$ws.Range("C31").Value = "percent"
$ws.Range("C32").Value = "plural"
$ws.Range("B33").Value = "Check the datetime data can be fetched"
$ws.Range("C33").Value = "datetime"
$ws.Range("B34").Value = "Check the currency data can be fetched"
$ws.Range("C34").Value = "currency"
$ws.Range("G34").Value = "US$201,703.54"
$ws.Range("G33").Value = "2017年11月20日 GMT+8 下午1:39:24"

$ws.Range("H31").Value = "CLDRJar"
$ws.Range("H32").Value = "CLDRJar"
$ws.Range("H33").Value = "CLDRJar"
$ws.Range("H34").Value = "CLDRJar"

$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32
$ws.Range("A34").Value = 33

$ws.Range("G30").Value = 201703.54199999999
$ws.Range("G30").NumberFormat = "#,##0.00"

$ws.Range("G31").Value = 0.23
$ws.Range("G31").NumberFormat = "0%"

# --- Restore view/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("A35").Select()
